$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2463960051175041
$ws.Range("C2").Value = 0.08281662434485339
$ws.Range("D2").Value = 366

$ws.Range("C3").Value = 0.4342047944680558

$ws.Range("B5").Value = [double]"3.127281843040124e-06"
$ws.Range("C5").Value = [double]"1.96923229072785e-06"
$ws.Range("D5").Value = 509

$ws.Range("B7").Value = [double]"-3.189883214351137e-24"
$ws.Range("C7").Value = [double]"-7.874071819255612e-24"
$ws.Range("D7").Value = 1967
$ws.Range("E7").Value = 1827

$ws.Range("B8").Value = 0.003589987327475219
$ws.Range("C8").Value = 0.002925460724059054
$ws.Range("D8").Value = 296
$ws.Range("E8").Value = 263

$ws.Range("B9").Value = 0.01671284838886192
$ws.Range("C9").Value = 0.01088913695132789
$ws.Range("D9").Value = 517
$ws.Range("E9").Value = 450

$ws.Range("B11").Value = [double]"7.263794386259762e-11"
$ws.Range("C11").Value = [double]"4.281413039908872e-11"
$ws.Range("D11").Value = 3705
$ws.Range("E11").Value = 3621

$ws.Range("B12").Value = 0.03105217670364126
$ws.Range("C12").Value = 0.03095296360307887

$ws.Range("B13").Value = 0.2725309561027044
$ws.Range("C13").Value = 0.2544924246032119
$ws.Range("D13").Value = 150
$ws.Range("E13").Value = 140

$ws.Range("B14").Value = 0.658629096420053
$ws.Range("C14").Value = 0.6329882972391327
$ws.Range("E14").Value = 48

$ws.Range("B16").Value = 0.4644080073711088
$ws.Range("C16").Value = 0.4253491238158417
$ws.Range("D16").Value = 107
$ws.Range("E16").Value = 93

$ws.Range("B17").Value = 0.618349052038167
$ws.Range("C17").Value = 0.6171135280926713
